# Normalizes "Tipo de Imóvel" (col A) text and "Tipo de Venda" (col O)
# labels, and fixes "Matrícula" (col Q) / "Inscrição Imobiliária" (col R)
# so they are stored as real numbers instead of text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: trim stray leading/trailing spaces around the property type ---
$lastRow = 25
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($null -ne $val) {
        $cell.Value = $val.ToString().Trim()
    }
}

# --- Column O: insert spaces into the concatenated "Tipo de Venda" labels ---
$oMap = @{
    2  = "Leilão SFI Caixa"
    3  = "Licitação Aberta Caixa"
    4  = "Licitação Aberta Caixa"
    5  = "Licitação Aberta Caixa"
    6  = "Licitação Aberta Caixa"
    7  = "Venda Online Caixa"
    8  = "Venda Online Caixa"
    9  = "Venda Online Caixa"
    10 = "Venda Online Caixa"
    11 = "Venda Direta Caixa"
    12 = "Venda Direta Caixa"
    13 = "Venda Direta"
    18 = "Leilão SFI Caixa"
    21 = "Leilão SFI Caixa"
    22 = "Leilão SFI Caixa"
    23 = "Leilão SFI Caixa"
    24 = "Venda Direta"
    25 = "Venda Direta"
}
foreach ($r in $oMap.Keys) {
    $ws.Cells.Item($r, 15).Value = $oMap[$r]
}

# --- Columns Q (Matrícula) and R (Inscrição Imobiliária): store as numbers ---
$qMap = @{
    2  = 16720
    3  = 5068
    4  = 8763
    5  = 16440
    6  = 14769
    7  = 5114
    8  = 5008
    9  = 15297
    10 = 14980
    11 = 10788
    12 = 9127
    18 = 13612
    21 = 3976
    22 = 13754
    23 = 13754
}
foreach ($r in $qMap.Keys) {
    $ws.Cells.Item($r, 17).Value = $qMap[$r]
}

$rMap = @{
    2  = 0
    7  = 220175
    8  = 221998
    9  = 222844
    10 = 224969
    11 = 3636062018400000
    12 = 219031
    18 = 0
    21 = 0
    22 = 0
    23 = 0
}
foreach ($r in $rMap.Keys) {
    $ws.Cells.Item($r, 18).Value = $rMap[$r]
}
